$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying weekly data rows (2-10) were refreshed: the set of observations
# got reordered (rotated) in place. Columns A,B,C,E,F,G,H,I,J stay identical for
# every row; only D (Fecha), K (Variedad), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion), R (Origen), S (Precio $/Kg) and
# T (Kg / unidad) change. Row 4 is untouched.

# Row 2
$ws.Cells.Item(2, 4).Value = [datetime]"2021-03-30"
$ws.Cells.Item(2, 11).Value = "Wonderfull"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 8
$ws.Cells.Item(2, 14).Value = 280000
$ws.Cells.Item(2, 15).Value = 300000
$ws.Cells.Item(2, 16).Value = 290000
$ws.Cells.Item(2, 17).Value = "`$/bins (400 kilos)"
$ws.Cells.Item(2, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(2, 19).Value = 725
$ws.Cells.Item(2, 20).Value = 400

# Row 3
$ws.Cells.Item(3, 4).Value = [datetime]"2021-04-21"
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 150
$ws.Cells.Item(3, 14).Value = 16000
$ws.Cells.Item(3, 15).Value = 18000
$ws.Cells.Item(3, 16).Value = 17000
$ws.Cells.Item(3, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(3, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(3, 19).Value = 1133
$ws.Cells.Item(3, 20).Value = 15

# Row 4 unchanged

# Row 5
$ws.Cells.Item(5, 4).Value = [datetime]"2021-05-18"
$ws.Cells.Item(5, 11).Value = "Wonderfull"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 16
$ws.Cells.Item(5, 14).Value = 240000
$ws.Cells.Item(5, 15).Value = 250000
$ws.Cells.Item(5, 16).Value = 245000
$ws.Cells.Item(5, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(5, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(5, 19).Value = 544
$ws.Cells.Item(5, 20).Value = 450

# Row 6
$ws.Cells.Item(6, 4).Value = [datetime]"2021-03-25"
$ws.Cells.Item(6, 11).Value = "Sin especificar"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 15
$ws.Cells.Item(6, 14).Value = 360000
$ws.Cells.Item(6, 15).Value = 360000
$ws.Cells.Item(6, 16).Value = 360000
$ws.Cells.Item(6, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(6, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(6, 19).Value = 800
$ws.Cells.Item(6, 20).Value = 450

# Row 7
$ws.Cells.Item(7, 4).Value = [datetime]"2021-03-11"
$ws.Cells.Item(7, 11).Value = "Wonderfull"
$ws.Cells.Item(7, 12).Value = "Segunda"
$ws.Cells.Item(7, 13).Value = 120
$ws.Cells.Item(7, 14).Value = 4800
$ws.Cells.Item(7, 15).Value = 4800
$ws.Cells.Item(7, 16).Value = 4800
$ws.Cells.Item(7, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(7, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(7, 19).Value = 1200
$ws.Cells.Item(7, 20).Value = 4

# Row 8
$ws.Cells.Item(8, 4).Value = [datetime]"2021-03-11"
$ws.Cells.Item(8, 11).Value = "Wonderfull"
$ws.Cells.Item(8, 12).Value = "Tercera"
$ws.Cells.Item(8, 13).Value = 80
$ws.Cells.Item(8, 14).Value = 4000
$ws.Cells.Item(8, 15).Value = 4000
$ws.Cells.Item(8, 16).Value = 4000
$ws.Cells.Item(8, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(8, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(8, 19).Value = 1000
$ws.Cells.Item(8, 20).Value = 4

# Row 9
$ws.Cells.Item(9, 4).Value = [datetime]"2021-04-26"
$ws.Cells.Item(9, 11).Value = "Wonderfull"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 24
$ws.Cells.Item(9, 14).Value = 220000
$ws.Cells.Item(9, 15).Value = 240000
$ws.Cells.Item(9, 16).Value = 230000
$ws.Cells.Item(9, 17).Value = "`$/bins (400 kilos)"
$ws.Cells.Item(9, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(9, 19).Value = 575
$ws.Cells.Item(9, 20).Value = 400

# Row 10
$ws.Cells.Item(10, 4).Value = [datetime]"2021-04-26"
$ws.Cells.Item(10, 11).Value = "Wonderfull"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 34
$ws.Cells.Item(10, 14).Value = 240000
$ws.Cells.Item(10, 15).Value = 240000
$ws.Cells.Item(10, 16).Value = 240000
$ws.Cells.Item(10, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(10, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(10, 19).Value = 533
$ws.Cells.Item(10, 20).Value = 450
